$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Add the four new 2019 semi-finals rows (205-208), reusing the same
# fill/style as the preceding data rows (A200:F204, style index 5) by
# copying that formatting down before writing the new values.
# ------------------------------------------------------------------
$ws.Range("A204:F204").Copy() | Out-Null
$ws.Range("A205:F208").PasteSpecial(-4122) | Out-Null

$ws.Range("A205").Value = 2019
$ws.Range("B205").Value = "Boston Bruins"
$ws.Range("C205").Value = "Columbus Blue Jackets"
$ws.Range("D205").Value = "eastern"
$ws.Range("E205").Value = "Boston Bruins"
$ws.Range("F205").Value = "semi-finals"

# Row 206 & 208 introduce new combined "Team/Team" shared strings - set
# the Team1 (col B) values first so the new shared-string entries land
# in the same order as the target workbook (46, 47, 48).
$ws.Range("A206").Value = 2019
$ws.Range("B206").Value = "Washington Capitals/Carolina Hurricanes"
$ws.Range("C206").Value = "New York Islanders"
$ws.Range("D206").Value = "eastern"
$ws.Range("F206").Value = "semi-finals"

$ws.Range("A207").Value = 2019
$ws.Range("B207").Value = "Dallas Stars"
$ws.Range("C207").Value = "St Louis Blues"
$ws.Range("D207").Value = "western"
$ws.Range("E207").Value = "St Louis Blues"
$ws.Range("F207").Value = "semi-finals"

$ws.Range("A208").Value = 2019
$ws.Range("B208").Value = "San Jose Sharks/Vegas Golden Knights"
$ws.Range("C208").Value = "Colorado Avalanche"
$ws.Range("D208").Value = "western"
$ws.Range("E208").Value = "San Jose Sharks/Vegas Golden Knights"
$ws.Range("F208").Value = "semi-finals"

$ws.Range("E206").Value = "New York Islanders/Washington Capitals"

# ------------------------------------------------------------------
# Widen Team1 (B) and Highest.Seed (E) columns so the longer combined
# team names fit - this splits the old merged B:C column-width entry.
# ------------------------------------------------------------------
$ws.Columns("B").ColumnWidth = 37.166666666666664
$ws.Columns("E").ColumnWidth = 36.666666666666664

# ------------------------------------------------------------------
# Update the view/selection to point at the new last row.
# ------------------------------------------------------------------
$ws.Range("E207").Select() | Out-Null
